# daily auto push: 2026-01-09 06:50 UTC
#
# The sheet is a flat "date / day-of-week / hour / rank" log, sorted
# chronologically. A new sample for 2026/01/09 (Friday) at hour 13 needs to
# be inserted right after the existing 2026/01/09 rows (596, 597), which
# pushes every following row down by one (598 -> 599, ..., 639 -> 640) and
# grows the used range from A1:D639 to A1:D640.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 598..639 down to 599..640, leaving row 598 blank.
$ws.Rows.Item(598).Insert()

# Fill the newly-opened row 598 with the inserted record.
# Column A holds dates formatted as plain text (e.g. "2026/12/29"), not
# real Excel dates, so force text interpretation with a leading apostrophe
# and then strip the auto-applied "Text" number format back to the sheet's
# normal (unstyled) look, matching every other data row.
$ws.Range("A598").Value = "'2026/01/09"
$ws.Range("A598").Style = "Normal"
$ws.Range("B598").Value = "金"
$ws.Range("C598").Value = 13
$ws.Range("D598").Value = 201
